$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B4 value from 0.7 to 0.8
$ws.Range("B4").Value = 0.8

# Row 2 height reverts to default (auto), removing the explicit 49.5 height
$ws.Rows("2").EntireRow.AutoFit()

# Update the selected cell to C4
$ws.Range("C4").Select()
